# Commit: 'matched styles for desserts in french docx'
# - Remove direct-formatting color overrides (0070C0) from the title and the two
#   Heading4 run/paragraph marks, now that the Heading3/Heading4 styles carry the
#   matching color themselves; split 'Les Cookies' into 'Les '+'Cookies' runs; drop
#   the stray empty paragraph under the title; re-point the heading styles' colors.
$d = $word.ActiveDocument

# --- Paragraph-level edits (process from the bottom of the doc upward so earlier
#     paragraph indices stay valid while later ones are still being edited) ---

# 'Préparation' Heading4 paragraph: drop the direct color formatting
$pPrep = $d.Paragraphs(12)
$xmlPrep = '<pkg:package xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage''><pkg:part pkg:name=''/word/document.xml'' pkg:contentType=''application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml''><pkg:xmlData><w:document xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:body><w:p><w:pPr><w:pStyle w:val="Heading4"/></w:pPr><w:r><w:t>Préparation</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pPrep.Range.InsertXML($xmlPrep)

# 'Ingrédients' Heading4 paragraph: drop the direct color formatting
$pIngr = $d.Paragraphs(3)
$xmlIngr = '<pkg:package xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage''><pkg:part pkg:name=''/word/document.xml'' pkg:contentType=''application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml''><pkg:xmlData><w:document xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:body><w:p><w:pPr><w:pStyle w:val="Heading4"/></w:pPr><w:r><w:t>Ingrédients</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pIngr.Range.InsertXML($xmlIngr)

# Empty paragraph right under the title: remove it entirely
$pEmpty = $d.Paragraphs(2)
$pEmpty.Range.Delete()

# Title ('Les Cookies « Classiques »', Heading3): drop direct color formatting and
# split 'Les Cookies' into two runs ('Les ' + 'Cookies') around the bookmark
$pTitle = $d.Paragraphs(1)
$xmlTitle = '<pkg:package xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage''><pkg:part pkg:name=''/word/document.xml'' pkg:contentType=''application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml''><pkg:xmlData><w:document xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:body><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:bookmarkStart w:id="0" w:name="_Toc395353017"/><w:r><w:t xml:space="preserve">Les </w:t></w:r><w:r><w:t>Cookies</w:t></w:r><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> « Classiques »</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pTitle.Range.InsertXML($xmlTitle)

# --- Style-level edits: re-point the heading styles (and their linked run styles)
#     at the new shared accent color 0070C0 ---
$wdColor0070C0 = 12611584  # RGB(0x00,0x70,0xC0) packed as BGR, per MS-OSHARED

$d.Styles("Heading 3").Font.Color = $wdColor0070C0
$d.Styles("Heading 3 Char").Font.Color = $wdColor0070C0
$d.Styles("Heading 4").Font.Color = $wdColor0070C0
$d.Styles("Heading 4 Char").Font.Color = $wdColor0070C0
